$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Units")
$ws.Activate()

$ws.Range("S1:V1").EntireColumn.Insert()

